$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shuakhevi")

# Add the 2023 column (K) to the right of the existing 2022 column (J),
# continuing the same header/data/formatting pattern as the rest of the table.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 700
$ws.Range("K5").Value = 78
$ws.Range("K6").Value = 622

# Copy formatting from column J (the previous last column) into K so the
# new column matches the existing look (borders, number format, font).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
